$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B24").Value = 6341
$ws.Range("D24").Value = 5939444
$ws.Range("E24").Value = 936.6730799558429
$ws.Range("F24").Value = 8.097511080804631
$ws.Range("H24").Value = 25.82043023599683
